$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1808', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1808', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1808', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1809', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1809', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1809', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1810', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1810', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1810', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1811', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1811', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1811', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1812', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1812', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1812', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1901', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1901', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1901', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1902', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1902', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1902', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1903', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1903', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1903', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1904', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1904', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1904', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1905', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1905', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1905', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1906', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1906', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1906', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1907', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1907', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1907', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1908', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1908', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1908', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1909', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1909', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1909', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1910', 40000, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1910', 40000, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1910', 40000, 1000000),
    @('1047371050', 'RONAL RAFAEL SALAS GUERRA', '1911', 22666, 1000000),
    @('1047489473', 'ELEAZAR DE JESUS ALCANTARA PEREZ', '1911', 22666, 1000000),
    @('13816726', 'DAGOBERTO COLEY ESTEVEZ', '1911', 22666, 1000000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
}
